$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows 5 and 6 (ACC-BU_A-0004, ACC-BU_A-0005), mirroring existing rows 2-4.
# Force the value_date column to text formatting before writing so the
# "2025-02-15" strings aren't auto-converted into date serial numbers,
# then restore the default "Normal" style so the cells end up styled the
# same as the other data rows (no explicit style index).
$ws.Range("E5:E6").NumberFormat = "@"

$ws.Range("A5").Value = 40500
$ws.Range("B5").Value = "BRL"
$ws.Range("C5").Value = "ACC-BU_A-0004"
$ws.Range("D5").Value = "Sample closure line 4 for BU_A"
$ws.Range("E5").Value = "2025-02-15"
$ws.Range("F5").Value = "BU_A"

$ws.Range("A6").Value = 50500
$ws.Range("B6").Value = "BRL"
$ws.Range("C6").Value = "ACC-BU_A-0005"
$ws.Range("D6").Value = "Sample closure line 5 for BU_A"
$ws.Range("E6").Value = "2025-02-15"
$ws.Range("F6").Value = "BU_A"

$ws.Range("E5:E6").Style = "Normal"
